$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dades_Període")

$ws.Range("H2").Value = "2026-02-20 10:15:29"
$ws.Range("H3").Value = "2026-02-20 10:15:30"
$ws.Range("H4").Value = "2026-02-20 10:15:30"
$ws.Range("H5").Value = "2026-02-20 10:15:30"
$ws.Range("H6").Value = "2026-02-20 10:15:30"
